# Auto-generated edit script: apply cell-level updates from Ifrit_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 104.55556
$ws.Range("I9").Value = 111.375
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 111.375
$ws.Range("L9").Value = 50
$ws.Range("M9").Value = 57.625
$ws.Range("N9").Value = -388
$ws.Range("H19").Value = 1662.9546
$ws.Range("I19").Value = 3074.7778
$ws.Range("J19").Value = 685.53845
$ws.Range("K19").Value = 3074.7778
$ws.Range("L19").Value = 685.53845
$ws.Range("M19").Value = -2899.7778
$ws.Range("N19").Value = -1035.53845
$ws.Range("H54").Value = 9800
$ws.Range("I54").Value = 9800
$ws.Range("K54").Value = 9800
$ws.Range("M54").Value = -9314
$ws.Range("H129").Value = 1231.4773
$ws.Range("I129").Value = 416
$ws.Range("J129").Value = 1313.025
$ws.Range("K129").Value = 1248
$ws.Range("L129").Value = 3939.075
$ws.Range("M129").Value = 3752
$ws.Range("N129").Value = -13939.075
$ws.Range("H132").Value = 252878.2
$ws.Range("I132").Value = 297339.97
$ws.Range("J132").Value = 928.1667
$ws.Range("K132").Value = 892019.9099999999
$ws.Range("L132").Value = 2784.5001
$ws.Range("M132").Value = -889489.9099999999
$ws.Range("N132").Value = -7844.5001
$ws.Range("H137").Value = 4746.9165
$ws.Range("I137").Value = 6450.1
$ws.Range("J137").Value = 2617.9375
$ws.Range("K137").Value = 19350.3
$ws.Range("L137").Value = 7853.8125
$ws.Range("M137").Value = -16800.3
$ws.Range("N137").Value = -12953.8125
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2503347.8
$ws.Range("I74").Value = 3333867.2
$ws.Range("J74").Value = 11789
$ws.Range("K74").Value = 3333867.2
$ws.Range("L74").Value = 11789
$ws.Range("M74").Value = -3332993.2
$ws.Range("N74").Value = -13537
$ws.Range("H77").Value = 2503347.8
$ws.Range("I77").Value = 3333867.2
$ws.Range("J77").Value = 11789
$ws.Range("K77").Value = 16669336
$ws.Range("L77").Value = 58945
$ws.Range("M77").Value = -16664968
$ws.Range("N77").Value = -67681
$ws.Range("H122").Value = 1449.1
$ws.Range("I122").Value = 1229.8462
$ws.Range("J122").Value = 1856.2858
$ws.Range("K122").Value = 3689.5386
$ws.Range("L122").Value = 5568.857400000001
$ws.Range("M122").Value = -1239.5386
$ws.Range("N122").Value = -10468.8574
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()
$ws.Range("H134").Value = 34930
$ws.Range("I134").Value = 43496.152
$ws.Range("J134").Value = 3112.8572
$ws.Range("K134").Value = 130488.456
$ws.Range("L134").Value = 9338.571599999999
$ws.Range("M134").Value = -127953.456
$ws.Range("N134").Value = -14408.5716
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1584.875
$ws.Range("J16").Value = 674.5
$ws.Range("L16").Value = 674.5
$ws.Range("N16").Value = -1248.5
$ws.Range("H113").Value = 1584.875
$ws.Range("J113").Value = 674.5
$ws.Range("L113").Value = 674.5
$ws.Range("N113").Value = -5014.5
$ws.Range("H131").Value = 40463.13
$ws.Range("J131").Value = 40463.13
$ws.Range("L131").Value = 40463.13
$ws.Range("N131").Value = -50543.13
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3258.8572
$ws.Range("I63").Value = 956
$ws.Range("J63").Value = 4180
$ws.Range("K63").Value = 2868
$ws.Range("L63").Value = 12540
$ws.Range("M63").Value = -2119
$ws.Range("N63").Value = -14038
$ws.Range("H64").Value = 1078470.1
$ws.Range("I64").Value = 982.4
$ws.Range("J64").Value = 1751900
$ws.Range("K64").Value = 2947.2
$ws.Range("L64").Value = 5255700
$ws.Range("M64").Value = -2677.2
$ws.Range("N64").Value = -5256240
$ws.Range("H66").Value = 3258.8572
$ws.Range("I66").Value = 956
$ws.Range("J66").Value = 4180
$ws.Range("K66").Value = 8604
$ws.Range("L66").Value = 37620
$ws.Range("M66").Value = -4860
$ws.Range("N66").Value = -45108
$ws.Range("H67").Value = 1078470.1
$ws.Range("I67").Value = 982.4
$ws.Range("J67").Value = 1751900
$ws.Range("K67").Value = 2947.2
$ws.Range("L67").Value = 5255700
$ws.Range("M67").Value = -2011.2
$ws.Range("N67").Value = -5257572
$ws.Range("H68").Value = 838.3605
$ws.Range("I68").Value = 616.1731
$ws.Range("J68").Value = 1178.1765
$ws.Range("K68").Value = 1848.5193
$ws.Range("L68").Value = 3534.5295
$ws.Range("M68").Value = -1037.5193
$ws.Range("N68").Value = -5156.529500000001
$ws.Range("H71").Value = 838.3605
$ws.Range("I71").Value = 616.1731
$ws.Range("J71").Value = 1178.1765
$ws.Range("K71").Value = 5545.5579
$ws.Range("L71").Value = 10603.5885
$ws.Range("M71").Value = -1489.5579
$ws.Range("N71").Value = -18715.5885
$ws.Range("H110").Value = 2903.8572
$ws.Range("J110").Value = 3660
$ws.Range("L110").Value = 10980
$ws.Range("N110").Value = -19160
$ws.Range("H131").Value = 1962937
$ws.Range("J131").Value = 2705278
$ws.Range("L131").Value = 8115834
$ws.Range("N131").Value = -8125914
$ws.Range("H133").Value = 7221.15
$ws.Range("I133").Value = 4930
$ws.Range("J133").Value = 7341.737
$ws.Range("K133").Value = 14790
$ws.Range("L133").Value = 22025.211
$ws.Range("M133").Value = -9730
$ws.Range("N133").Value = -32145.211
$ws.Range("H139").Value = 2691.139
$ws.Range("I139").Value = 1403.6818
$ws.Range("J139").Value = 4714.2856
$ws.Range("K139").Value = 4211.0454
$ws.Range("L139").Value = 14142.8568
$ws.Range("M139").Value = 928.9546
$ws.Range("N139").Value = -24422.8568
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4255
$ws.Range("I122").Value = 4881.0527
$ws.Range("J122").Value = 2555.7144
$ws.Range("K122").Value = 14643.1581
$ws.Range("L122").Value = 7667.1432
$ws.Range("M122").Value = -12193.1581
$ws.Range("N122").Value = -12567.1432
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 37287.5
$ws.Range("I75").Value = 3990.5
$ws.Range("J75").Value = 48386.5
$ws.Range("K75").Value = 3990.5
$ws.Range("L75").Value = 48386.5
$ws.Range("M75").Value = -3054.5
$ws.Range("N75").Value = -50258.5
$ws.Range("H78").Value = 37287.5
$ws.Range("I78").Value = 3990.5
$ws.Range("J78").Value = 48386.5
$ws.Range("K78").Value = 11971.5
$ws.Range("L78").Value = 145159.5
$ws.Range("M78").Value = -7291.5
$ws.Range("N78").Value = -154519.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 718577.9
$ws.Range("I52").Value = 2500750
$ws.Range("J52").Value = 5709
$ws.Range("K52").Value = 2500750
$ws.Range("L52").Value = 5709
$ws.Range("M52").Value = -2500524
$ws.Range("N52").Value = -6161
$ws.Range("H132").Value = 3744.803
$ws.Range("I132").Value = 4174.4287
$ws.Range("J132").Value = 1338.9
$ws.Range("K132").Value = 12523.2861
$ws.Range("L132").Value = 4016.7
$ws.Range("M132").Value = -9993.286100000001
$ws.Range("N132").Value = -9076.700000000001
$ws.Range("H136").Value = 68037.61
$ws.Range("I136").Value = 15548.429
$ws.Range("J136").Value = 251749.75
$ws.Range("K136").Value = 46645.287
$ws.Range("L136").Value = 755249.25
$ws.Range("M136").Value = -44095.287
$ws.Range("N136").Value = -760349.25
